# Remove the liquidacion (CSG) columns BH:BV for rows 6 and 7.
# These columns were populated with per-row liquidation data (BH = CSG date
# string, BI..BV = numeric computations). The fix removes this data for
# shipment rows that have no CSG / liquidacion items, leaving only the
# existing BG (liquidation status) column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BH6:BV6").ClearContents()
$ws.Range("BH7:BV7").ClearContents()
